$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "E"=2; "F"=0.6666666666666666; "G"=0.1636203333333333; "H"=0.490861; "I"=0.3345941539187231; "J"=0.3345941539187231; "M"=1.918906333333333; "N"=5.756718999999999; "O"=0.006524019162508824; "P"=0.006524019162508824; "Q"=0.3139720938954444; "R"=2.825748845059; "S"=0.002182898671829176; "T"=0.002182898671829176 }
    3 = @{ "E"=2; "F"=0.6666666666666666; "G"=0.1636203333333333; "H"=0.490861; "I"=0.3345941539187231; "J"=0.3345941539187231; "O"=0.6163557430885885; "P"=0.6163557430885885; "Q"=29.66246701942378; "R"=266.962203174814; "S"=0.2062290283716721; "T"=0.2062290283716721 }
    4 = @{ "E"=2; "F"=0.6666666666666666; "G"=0.1636203333333333; "H"=0.490861; "I"=0.3345941539187231; "J"=0.3345941539187231; "M"=29.04767233333333; "N"=87.143017; "O"=0.09875811426384234; "P"=0.09875811426384236; "Q"=4.752789829737445; "R"=42.775108467637; "S"=0.03304388768471891; "T"=0.03304388768471891 }
    5 = @{ "E"=2; "F"=0.6666666666666666; "G"=0.1636203333333333; "H"=0.490861; "I"=0.3345941539187231; "J"=0.3345941539187231; "M"=81.87450533333333; "N"=245.623516; "O"=0.2783621234850603; "P"=0.2783621234850603; "Q"=13.39633385414178; "R"=120.567004687276; "S"=0.09313833919050286; "T"=0.09313833919050286 }
    6 = @{ "E"=2; "F"=0.6666666666666666; "G"=0.325391; "H"=0.9761730000000001; "I"=0.665405846081277; "J"=0.665405846081277; "M"=1.918906333333333; "N"=5.756718999999999; "O"=0.006524019162508824; "P"=0.006524019162508824; "Q"=0.6243948507096667; "R"=5.619553656387; "S"=0.004341120490679647; "T"=0.004341120490679647 }
    7 = @{ "E"=2; "F"=0.6666666666666666; "G"=0.325391; "H"=0.9761730000000001; "I"=0.665405846081277; "J"=0.665405846081277; "O"=0.6163557430885885; "P"=0.6163557430885885; "Q"=58.98961094434468; "R"=530.9064984991021; "S"=0.4101267147169164; "T"=0.4101267147169164 }
    8 = @{ "E"=2; "F"=0.6666666666666666; "G"=0.325391; "H"=0.9761730000000001; "I"=0.665405846081277; "J"=0.665405846081277; "M"=29.04767233333333; "N"=87.143017; "O"=0.09875811426384234; "P"=0.09875811426384236; "Q"=9.451851148215669; "R"=85.066660333941; "S"=0.06571422657912344; "T"=0.06571422657912344 }
    9 = @{ "E"=2; "F"=0.6666666666666666; "G"=0.325391; "H"=0.9761730000000001; "I"=0.665405846081277; "J"=0.665405846081277; "M"=81.87450533333333; "N"=245.623516; "O"=0.2783621234850603; "P"=0.2783621234850603; "Q"=26.64122716491867; "R"=239.771044484268; "S"=0.1852237842945574; "T"=0.1852237842945574 }
}

foreach ($rowKey in $updates.Keys) {
    $rowData = $updates[$rowKey]
    foreach ($colKey in $rowData.Keys) {
        $ws.Range("$colKey$rowKey").Value = $rowData[$colKey]
    }
}
